$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.005.31"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.829.89"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").Value = "'0.9984"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'243.90"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "'0.6318"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("D7").Value = "'0.9992"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.07521"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "'22.87"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").Value = "'0.07727"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "1.841.64"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "'4.995"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "'0.6709"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "'0.000009699"
$ws.Range("E16").Value = "  +7.29%  "
$ws.Range("D17").Value = "'6.078"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "29.024.86"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "'226.66"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "'0.9984"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'7.175"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "'0.9991"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'159.75"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'0.1406"
$ws.Range("E25").Value = "  +3.49%  "
$ws.Range("D26").Value = "'8.536"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'1.495"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'4.121"
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").Value = "'4.073"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").Value = "'1.200"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").Value = "'0.05366"
$ws.Range("E32").Value = "  +2.94%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").Value = "'1.139"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").Value = "'2.655"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Value = "1.245.05"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.752"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01785"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").Value = "'6.595"
$ws.Range("E40").Value = "  +3.62%  "
$ws.Range("D41").Value = "'0.9030"
$ws.Range("E41").Value = "  +1.42%  "
$ws.Range("D42").Value = "'0.9986"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'101.51"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "1.982.20"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'64.86"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000122"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("D47").Value = "'0.5101"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("D49").Value = "'9.035"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").Value = "'0.05765"
$ws.Range("E51").Value = "  +0.22%  "
